$wb = $excel.ActiveWorkbook

# sum_response_time_Results: recomputed REST statistics after removing
# the REST optimization (use cases "3" and "5" affected)
$ws1 = $wb.Worksheets.Item("sum_response_time_Results")

# Use Case "3" row
$ws1.Range("C3").Value = 18.76
$ws1.Range("D3").Value = 3.597554725087584
$ws1.Range("G3").Value = 0.001163470059104413
$ws1.Range("K3").Value = 3183139353988279 * [Math]::Pow(10, -32)

# Use Case "5" row
$ws1.Range("C5").Value = 16.9
$ws1.Range("D5").Value = 6.974955196988723
$ws1.Range("G5").Value = 1980423693536049 * [Math]::Pow(10, -24)
$ws1.Range("K5").Value = 0.0105601122249601

# total_data_transferred_Results: updated REST transfer totals
$ws2 = $wb.Worksheets.Item("total_data_transferred_Results")

$ws2.Range("C3").Value = 9551
$ws2.Range("C5").Value = 1201
